$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-66)
# from serial date 45185 (2023-09-16) to 45204 (2023-10-05).
$ws.Range("C2:C66").Value = 45204
